$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) "总计" (sheet1): shift the existing 6 data rows down by one row and
#    insert a brand-new top row for 2022-Q3.
# ---------------------------------------------------------------------------
$totalSheet = $wb.Worksheets.Item("总计")

$totalSheet.Range("A7:D7").Copy($totalSheet.Range("A8:D8"))
$totalSheet.Range("A6:D6").Copy($totalSheet.Range("A7:D7"))
$totalSheet.Range("A5:D5").Copy($totalSheet.Range("A6:D6"))
$totalSheet.Range("A4:D4").Copy($totalSheet.Range("A5:D5"))
$totalSheet.Range("A3:D3").Copy($totalSheet.Range("A4:D4"))
$totalSheet.Range("A2:D2").Copy($totalSheet.Range("A3:D3"))

$totalSheet.Range("A2").Value = 0
$totalSheet.Range("B2").Value = "2022-Q3"
$totalSheet.Range("C2").Value = 5
$totalSheet.Range("D2").Value = 0.05

$totalSheet.Range("A3").Value = 1
$totalSheet.Range("A4").Value = 2
$totalSheet.Range("A5").Value = 3
$totalSheet.Range("A6").Value = 4
$totalSheet.Range("A7").Value = 5
$totalSheet.Range("A8").Value = 6

# ---------------------------------------------------------------------------
# 2) Insert a brand-new "2022-Q3" sheet right before "2022-Q2" by duplicating
#    the "2022-Q2" sheet (so headers/styles match exactly), then overwrite
#    its contents with the 2022-Q3 fund-holding data.
# ---------------------------------------------------------------------------
$q2Sheet = $wb.Worksheets.Item("2022-Q2")
$q2Sheet.Copy($q2Sheet)

$q3Sheet = $wb.Worksheets.Item(2)
$q3Sheet.Name = "2022-Q3"

# grow the copied template from 4 data rows to 5 data rows
$q3Sheet.Range("A5:H5").Copy($q3Sheet.Range("A6:H6"))

# Row 2
$q3Sheet.Range("A2").Value = 0
$q3Sheet.Range("B2").NumberFormat = "@"
$q3Sheet.Range("B2").Value = "002152"
$q3Sheet.Range("C2").NumberFormat = "@"
$q3Sheet.Range("C2").Value = "华宝核心优势灵活配置混合A"
$q3Sheet.Range("D2").NumberFormat = "@"
$q3Sheet.Range("D2").Value = "0.47"
$q3Sheet.Range("E2").NumberFormat = "@"
$q3Sheet.Range("E2").Value = "91.02"
$q3Sheet.Range("F2").NumberFormat = "@"
$q3Sheet.Range("F2").Value = "7.60"
$q3Sheet.Range("G2").NumberFormat = "@"
$q3Sheet.Range("G2").Value = "0.0357"
$q3Sheet.Range("H2").Value = 1

# Row 3
$q3Sheet.Range("A3").Value = 1
$q3Sheet.Range("B3").NumberFormat = "@"
$q3Sheet.Range("B3").Value = "005360"
$q3Sheet.Range("C3").NumberFormat = "@"
$q3Sheet.Range("C3").Value = "汇安资产轮动灵活配置混合"
$q3Sheet.Range("D3").NumberFormat = "@"
$q3Sheet.Range("D3").Value = "0.12"
$q3Sheet.Range("E3").NumberFormat = "@"
$q3Sheet.Range("E3").Value = "80.44"
$q3Sheet.Range("F3").NumberFormat = "@"
$q3Sheet.Range("F3").Value = "6.20"
$q3Sheet.Range("G3").NumberFormat = "@"
$q3Sheet.Range("G3").Value = "0.0074"
$q3Sheet.Range("H3").Value = 3

# Row 4
$q3Sheet.Range("A4").Value = 2
$q3Sheet.Range("B4").NumberFormat = "@"
$q3Sheet.Range("B4").Value = "006231"
$q3Sheet.Range("C4").NumberFormat = "@"
$q3Sheet.Range("C4").Value = "国融融君灵活配置混合A"
$q3Sheet.Range("D4").NumberFormat = "@"
$q3Sheet.Range("D4").Value = "0.10"
$q3Sheet.Range("E4").NumberFormat = "@"
$q3Sheet.Range("E4").Value = "55.44"
$q3Sheet.Range("F4").NumberFormat = "@"
$q3Sheet.Range("F4").Value = "2.17"
$q3Sheet.Range("G4").NumberFormat = "@"
$q3Sheet.Range("G4").Value = "0.0022"
$q3Sheet.Range("H4").Value = 5

# Row 5
$q3Sheet.Range("A5").Value = 3
$q3Sheet.Range("B5").NumberFormat = "@"
$q3Sheet.Range("B5").Value = "016461"
$q3Sheet.Range("C5").NumberFormat = "@"
$q3Sheet.Range("C5").Value = "华宝核心优势灵活配置混合C"
$q3Sheet.Range("D5").NumberFormat = "@"
$q3Sheet.Range("D5").Value = "0.00"
$q3Sheet.Range("E5").NumberFormat = "@"
$q3Sheet.Range("E5").Value = "91.02"
$q3Sheet.Range("F5").NumberFormat = "@"
$q3Sheet.Range("F5").Value = "7.60"
$q3Sheet.Range("G5").Value = 0
$q3Sheet.Range("H5").Value = 1

# Row 6
$q3Sheet.Range("A6").Value = 4
$q3Sheet.Range("B6").NumberFormat = "@"
$q3Sheet.Range("B6").Value = "006232"
$q3Sheet.Range("C6").NumberFormat = "@"
$q3Sheet.Range("C6").Value = "国融融君灵活配置混合C"
$q3Sheet.Range("D6").NumberFormat = "@"
$q3Sheet.Range("D6").Value = "0.00"
$q3Sheet.Range("E6").NumberFormat = "@"
$q3Sheet.Range("E6").Value = "55.44"
$q3Sheet.Range("F6").NumberFormat = "@"
$q3Sheet.Range("F6").Value = "2.17"
$q3Sheet.Range("G6").Value = 0
$q3Sheet.Range("H6").Value = 5
